$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.517.75'
$ws.Range("E2").Value = '  -0.06%  '

$ws.Range("D3").Value = '1.913.44'
$ws.Range("E3").Value = '  +0.26%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").Value = '0.705'
$ws.Range("E5").Value = '  +6.24%  '

$ws.Range("D6").Value = '247.14'
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").Value = '40.68'

$ws.Range("D9").Value = '0.356'
$ws.Range("E9").Value = '  +3.58%  '

$ws.Range("D10").Value = '52.68'
$ws.Range("E10").Value = '  +7.22%  '

$ws.Range("E11").Value = '  +2.22%  '

$ws.Range("D12").Value = '0.0991'
$ws.Range("E12").Value = '  -1.02%  '

$ws.Range("D13").Value = '2.189.42'
$ws.Range("E13").Value = '  +0.24%  '

$ws.Range("D14").Value = '12.71'
$ws.Range("E14").Value = '  +2.73%  '

$ws.Range("D15").Value = '0.717'
$ws.Range("E15").Value = '  +2.24%  '

$ws.Range("D16").Value = '1.910.66'
$ws.Range("E16").Value = '  -0.25%  '

$ws.Range("E17").Value = '  +0.94%  '

$ws.Range("D18").Value = '35.503.78'
$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("D19").Value = '73.28'
$ws.Range("E19").Value = '  +1.23%  '

$ws.Range("D20").Value = '0.0₃0827'
$ws.Range("E20").Value = '  -0.61%  '

$ws.Range("D21").Value = '13.13'
$ws.Range("E21").Value = '  +3.83%  '

$ws.Range("D22").Value = '242.45'
$ws.Range("E22").Value = '  -0.61%  '

$ws.Range("D23").Value = '5.07'
$ws.Range("E23").Value = '  +4.32%  '

$ws.Range("E24").Value = '  -0.26%  '

$ws.Range("E25").Value = '  +1.22%  '

$ws.Range("E26").Value = '  +4.58%  '

$ws.Range("D27").Value = '168.99'
$ws.Range("E27").Value = '  -1.53%  '

$ws.Range("D28").Value = '8.64'
$ws.Range("E28").Value = '  +1.65%  '

$ws.Range("D29").Value = '18.78'
$ws.Range("E29").Value = '  +2.63%  '

$ws.Range("E30").Value = '  +2.26%  '

$ws.Range("D31").Value = '4.186.74'
$ws.Range("E31").Value = '  +21.13%  '

$ws.Range("E32").Value = '  +1.65%  '

$ws.Range("D33").Value = '0.0576'
$ws.Range("E33").Value = '  +0.87%  '

$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.90'
$ws.Range("E34").Value = '  +10.21%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = '4.21'
$ws.Range("E35").Value = '  -0.66%  '

$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("D37").Value = '0.911'
$ws.Range("E37").Value = '  -5.77%  '

$ws.Range("E38").Value = '  +11.41%  '

$ws.Range("E39").Value = '  +0.71%  '

$ws.Range("D40").Value = '17.40'
$ws.Range("E40").Value = '  +10.94%  '

$ws.Range("D41").Value = '98.48'
$ws.Range("E41").Value = '  +6.59%  '

$ws.Range("E42").Value = '  +2.92%  '

$ws.Range("E43").Value = '  +2.20%  '

$ws.Range("E44").Value = '  +1.85%  '

$ws.Range("D45").Value = '1.353.67'
$ws.Range("E45").Value = '  +0.38%  '

$ws.Range("D46").Value = '2.45'
$ws.Range("E46").Value = '  +2.39%  '

$ws.Range("E47").Value = '  +0.20%  '

$ws.Range("E48").Value = '  +0.74%  '

$ws.Range("D49").Value = '45.80'
$ws.Range("E49").Value = '  -3.89%  '

$ws.Range("D50").Value = '12.29'
$ws.Range("E50").Value = '  -2.61%  '

$ws.Range("E51").Value = '  -0.44%  '
